$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ten_lists")

# Add new "start DD" / "start SD" values into column E and K for each of the
# five list blocks (header rows 3, 10, 17, 24, 31).
$ws.Range("E3").Value  = "start DD"
$ws.Range("K3").Value  = "start SD"

$ws.Range("E10").Value = "start SD"
$ws.Range("K10").Value = "start DD"

$ws.Range("E17").Value = "start DD"
$ws.Range("K17").Value = "start SD"

$ws.Range("E24").Value = "start SD"
$ws.Range("K24").Value = "start DD"

$ws.Range("E31").Value = "start DD"
$ws.Range("K31").Value = "start DD"

# Match the author's final selection as recorded in the saved workbook.
$ws.Range("K31").Select()
